$d = $word.ActiveDocument

# Locate the end of the paragraph that ends with
# "...即包的抽象程度跟它的稳定性成正比。" (the end of section "4. 稳定抽象原则").
# Find.Execute collapses the range it's called on to the matched text span,
# so $rng.End lands right after the last character of the match and before
# its paragraph mark.
$rng = $d.Content
$found = $rng.Find.Execute("即包的抽象程度跟它的稳定性成正比。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor paragraph for the Stable Dependencies Principle insertion."
}
$insertPos = $rng.End

# Build the new content as an OOXML fragment (via InsertXML) so the inserted
# paragraphs/runs come out with the exact structure used elsewhere in this
# document (an empty separator paragraph, a "5. <title>" heading paragraph
# made of two runs, and a body paragraph), matching the authoring pattern of
# the preceding "4. 稳定抽象原则" section.
$ins = $d.Range($insertPos, $insertPos)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
         '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
             '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
               '<w:body>' + `
                 '<w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>' + `
                 '<w:p>' + `
                   '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
                   '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">5. </w:t></w:r>' + `
                   '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>稳定依赖原则</w:t></w:r>' + `
                 '</w:p>' + `
                 '<w:p>' + `
                   '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
                   '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>包之间的依赖关系都应该是稳定方向依赖的，包要依赖的包要比自己更具有稳定性。</w:t></w:r>' + `
                 '</w:p>' + `
               '</w:body>' + `
             '</w:document>' + `
           '</pkg:xmlData>' + `
         '</pkg:part>' + `
       '</pkg:package>'

$null = $ins.InsertXML($xml)

Write-Output "Inserted Stable Dependencies Principle section."
